$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit shuffles the text labels in column A between several rows
# (the counts in column B stay with their original row). Apply the new
# label for each row that changed.
$ws.Range("A19").Value = "крамными товар"
$ws.Range("A20").Value = "небогатый товар"
$ws.Range("A24").Value = "нужный товар"
$ws.Range("A25").Value = "пушной товар"
$ws.Range("A27").Value = "набойчатый товар"
$ws.Range("A28").Value = "медный товар"
$ws.Range("A29").Value = "суровский товар"
$ws.Range("A30").Value = "питейный припасы"
$ws.Range("A32").Value = "недорогой товар"
$ws.Range("A36").Value = "произрастание"
$ws.Range("A37").Value = "купецкий товар"
$ws.Range("A38").Value = "заморский товар"
$ws.Range("A39").Value = "меховой товар"
$ws.Range("A40").Value = "харчевой припасы"
$ws.Range("A42").Value = "рукодельный товар"
$ws.Range("A43").Value = "домовый товар"
